# Updated cryptos list (prices + 1h volume deltas) and re-ranked three
# coins (Aave / Mantle / BabyDogeCoin) that swapped rank order, per the
# GitHub Actions refresh job commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.230.81'
$ws.Range('E2').Value = '  -0.13%  '
# Row 3
$ws.Range('D3').Value = '1.844.16'
$ws.Range('E3').Value = '  +0.04%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('E4').Value = '  +0.02%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.91'
$ws.Range('E5').Value = '  +0.92%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6644'
$ws.Range('E6').Value = '  -0.85%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.0000'
$ws.Range('E7').Value = '  +0.04%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.94'
$ws.Range('E8').Value = '  +6.57%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07471'
$ws.Range('E9').Value = '  +0.41%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2965'
$ws.Range('E10').Value = '  +0.03%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.42'
$ws.Range('E11').Value = '  +2.88%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07756'
$ws.Range('E12').Value = '  +0.38%  '
# Row 13
$ws.Range('D13').Value = '1.851.05'
$ws.Range('E13').Value = '  +2.07%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.032'
$ws.Range('E14').Value = '  +0.16%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6760'
$ws.Range('E15').Value = '  -0.42%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.42'
$ws.Range('E16').Value = '  -3.45%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.182'
$ws.Range('E17').Value = '  +0.05%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008756'
$ws.Range('E18').Value = '  +5.42%  '
# Row 19
$ws.Range('D19').Value = '29.200.49'
$ws.Range('E19').Value = '  +0.31%  '
# Row 20
$ws.Range('D20').Value = '2.096.57'
$ws.Range('E20').Value = '  +3.00%  '
# Row 21
$ws.Range('E21').Value = '  +0.07%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.57'
$ws.Range('E22').Value = '  -0.61%  '
# Row 23
$ws.Range('E23').Value = '  +0.14%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.214'
$ws.Range('E24').Value = '  +0.11%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.0000'
$ws.Range('E25').Value = '  +0.02%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.42'
$ws.Range('E26').Value = '  -0.94%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.645'
$ws.Range('E27').Value = '  -0.62%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1403'
$ws.Range('E28').Value = '  +0.39%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.09'
$ws.Range('E29').Value = '  +0.26%  '
# Row 30
$ws.Range('E30').Value = '  +0.01%  '
# Row 31
$ws.Range('E31').Value = '  -1.21%  '
# Row 32
$ws.Range('E32').Value = '  -0.83%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.193'
$ws.Range('E33').Value = '  +0.04%  '
# Row 34
$ws.Range('E34').Value = '  +1.00%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.855'
$ws.Range('E35').Value = '  -2.02%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7480'
$ws.Range('E36').Value = '  -1.37%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.162'
$ws.Range('E37').Value = '  +1.38%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.643'
$ws.Range('E38').Value = '  -1.41%  '
# Row 39
$ws.Range('D39').Value = '1.303.12'
$ws.Range('E39').Value = '  -2.22%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01798'
$ws.Range('E40').Value = '  -0.67%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.759'
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.423'
$ws.Range('E42').Value = '  +8.14%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9073'
$ws.Range('E43').Value = '  -2.06%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9992'
$ws.Range('E44').Value = '  -0.22%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.71'
$ws.Range('E45').Value = '  +0.30%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.08235'
$ws.Range('E46').Value = '  +2.89%  '
# Row 47
$ws.Range('D47').Value = '1.996.50'
$ws.Range('E47').Value = '  +2.93%  '
# Row 48
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.38'
$ws.Range('E48').Value = '  +2.63%  '
# Row 49
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5142'
$ws.Range('E49').Value = '  -0.35%  '
# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000121'
$ws.Range('E50').Value = '  -1.30%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.756'
$ws.Range('E51').Value = '  -0.86%  '
